# Apply the row-level data corrections described by the commit diff.
# (Rows 5/6, 7/9, 8/14 swap their record content; rows 10-13 rotate;
# "Taxonsorteringsordning" (col B) for every record increases by 14.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 90808

# Row 3
$ws.Range("B3").Value = 89794

# Row 4
$ws.Range("B4").Value = 91002

# Row 5
$ws.Range("A5").Value = 112445391
$ws.Range("B5").Value = 89478
$ws.Range("D5").Value = 'NT'
$ws.Range("E5").Value = 3242
$ws.Range("F5").Value = 'Vitplätt'
$ws.Range("G5").Value = 'Chaetodermella luna'
$ws.Range("H5").Value = '(Romell ex D.P.Rogers & H.S.Jacks.) Rauschert'
$ws.Range("Q5").Value = 749781
$ws.Range("R5").Value = 7536174
$ws.Range("S5").Value = 1
$ws.Range("Z5").Value = '14:51'
$ws.Range("AB5").Value = '14:51'
$ws.Range("AW5").Value = 'Stefan Andersson'
$ws.Range("AX5").Value = 'Stefan Andersson, Christina Boyd, per-erik mukka'

# Row 6
$ws.Range("A6").Value = 112446488
$ws.Range("B6").Value = 89517
$ws.Range("D6").Value = 'LC'
$ws.Range("E6").Value = 5447
$ws.Range("F6").Value = 'Vedticka'
$ws.Range("G6").Value = 'Fuscoporia viticola'
$ws.Range("H6").Value = '(Schwein.) Murrill'
$ws.Range("Q6").Value = 749792
$ws.Range("R6").Value = 7536087
$ws.Range("S6").Value = 10
$ws.Range("Z6").Value = $null
$ws.Range("AB6").Value = $null
$ws.Range("AW6").Value = 'per-erik mukka'
$ws.Range("AX6").Value = 'per-erik mukka, Christina Boyd, Stefan Andersson'

# Row 7
$ws.Range("A7").Value = 112443040
$ws.Range("B7").Value = 90808
$ws.Range("Q7").Value = 749867
$ws.Range("R7").Value = 7536102
$ws.Range("S7").Value = 10
$ws.Range("Z7").Value = '13:48'
$ws.Range("AB7").Value = '13:48'

# Row 8
$ws.Range("A8").Value = 112439340
$ws.Range("B8").Value = 95707
$ws.Range("E8").Value = 221941
$ws.Range("F8").Value = 'Plattlummer'
$ws.Range("G8").Value = 'Lycopodium complanatum'
$ws.Range("H8").Value = 'L.'
$ws.Range("P8").Value = 'Rautusakaravägen, T lm'
$ws.Range("Q8").Value = 749873
$ws.Range("R8").Value = 7536596
$ws.Range("S8").Value = 10
$ws.Range("Z8").Value = '12:14'
$ws.Range("AB8").Value = '12:14'

# Row 9
$ws.Range("A9").Value = 112444020
$ws.Range("B9").Value = 90808
$ws.Range("Q9").Value = 749749
$ws.Range("R9").Value = 7536225
$ws.Range("S9").Value = 50
$ws.Range("Z9").Value = '14:18'
$ws.Range("AB9").Value = '14:18'

# Row 10
$ws.Range("A10").Value = 112445270
$ws.Range("B10").Value = 90826
$ws.Range("D10").Value = 'LC'
$ws.Range("E10").Value = 4366
$ws.Range("F10").Value = 'Skarp dropptaggsvamp'
$ws.Range("G10").Value = 'Hydnellum peckii'
$ws.Range("H10").Value = 'Banker'
$ws.Range("Q10").Value = 749804
$ws.Range("R10").Value = 7536188
$ws.Range("S10").Value = 25
$ws.Range("Z10").Value = '14:48'
$ws.Range("AB10").Value = '14:48'

# Row 11
$ws.Range("A11").Value = 112446507
$ws.Range("B11").Value = 89478
$ws.Range("D11").Value = 'NT'
$ws.Range("E11").Value = 3242
$ws.Range("F11").Value = 'Vitplätt'
$ws.Range("G11").Value = 'Chaetodermella luna'
$ws.Range("H11").Value = '(Romell ex D.P.Rogers & H.S.Jacks.) Rauschert'
$ws.Range("Q11").Value = 749840
$ws.Range("S11").Value = 100
$ws.Range("Z11").Value = '15:25'
$ws.Range("AB11").Value = '15:25'

# Row 12
$ws.Range("A12").Value = 112444819
$ws.Range("B12").Value = 90806
$ws.Range("D12").Value = 'NT'
$ws.Range("E12").Value = 4361
$ws.Range("F12").Value = 'Orange taggsvamp'
$ws.Range("G12").Value = 'Hydnellum aurantiacum'
$ws.Range("H12").Value = '(Batsch:Fr.) P.Karst.'
$ws.Range("Q12").Value = 749749
$ws.Range("R12").Value = 7536225
$ws.Range("S12").Value = 50
$ws.Range("Z12").Value = '14:37'
$ws.Range("AB12").Value = '14:37'

# Row 13
$ws.Range("A13").Value = 112437506
$ws.Range("B13").Value = 89794
$ws.Range("D13").Value = 'VU'
$ws.Range("E13").Value = 65
$ws.Range("F13").Value = 'Fläckporing'
$ws.Range("G13").Value = 'Anthoporia albobrunnea'
$ws.Range("H13").Value = '(Romell) Karasiński & Niemelä'
$ws.Range("Q13").Value = 749872
$ws.Range("R13").Value = 7536255
$ws.Range("S13").Value = 25
$ws.Range("Z13").Value = '11:27'
$ws.Range("AB13").Value = '11:27'

# Row 14
$ws.Range("A14").Value = 112444216
$ws.Range("B14").Value = 90814
$ws.Range("E14").Value = 4364
$ws.Range("F14").Value = 'Dropptaggsvamp'
$ws.Range("G14").Value = 'Hydnellum ferrugineum'
$ws.Range("H14").Value = '(Fr.:Fr.) P. Karst.'
$ws.Range("P14").Value = 'Paurankivaravägen, T lm'
$ws.Range("Q14").Value = 749749
$ws.Range("R14").Value = 7536225
$ws.Range("S14").Value = 25
$ws.Range("Z14").Value = '14:24'
$ws.Range("AB14").Value = '14:24'
